$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "unit" shared strings: "IN" -> "inch", "in" -> "INch".
# Cells currently showing "in" (F14) should end up reading "inch",
# and cells currently showing "IN" (F16) should end up reading "INch".
$ws.Range("F14").Value = "inch"
$ws.Range("F16").Value = "INch"

# Update the selected cell on the sheet view from E17 to F17.
$ws.Range("F17").Select()
